$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)  # Column D = Speaker
    $val = $cell.Text
    if ($val -eq "Davis") {
        $cell.Value = "T"
    }
    elseif ($val -eq "Student") {
        $cell.Value = "S"
    }
}
